# Update odds values on Sheet1 (ActiveSheet) to match the 2024-11-01 refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Central Cordoba vs San Lorenzo)
$ws.Range("G2").Value = 2.25
$ws.Range("H2").Value = 2.88
$ws.Range("I2").Value = 3.8
$ws.Range("J2").Value = 3.25
$ws.Range("K2").Value = 1.73
$ws.Range("M2").Value = 1.17
$ws.Range("N2").Value = 5
$ws.Range("X2").Value = 8.5
$ws.Range("AA2").Value = 29
$ws.Range("AG2").Value = 6.5
$ws.Range("AJ2").Value = 41
$ws.Range("AN2").Value = 4
$ws.Range("AP2").Value = 41
$ws.Range("AR2").Value = 126
$ws.Range("AT2").Value = 1.91
$ws.Range("AZ2").Value = 101

# Row 5 (Argentinos Jrs vs Velez Sarsfield)
$ws.Range("M5").Value = 1.08
$ws.Range("N5").Value = 7.5
$ws.Range("X5").Value = 15
$ws.Range("AA5").Value = 29
$ws.Range("AC5").Value = 7.5
$ws.Range("AS5").Value = 251

# Row 7 (Fluminense vs Gremio)
$ws.Range("G7").Value = 2.2
$ws.Range("I7").Value = 3.5
$ws.Range("U7").Value = 2.05
$ws.Range("V7").Value = 1.7
$ws.Range("Z7").Value = 21
$ws.Range("AH7").Value = 15
$ws.Range("BA7").Value = 101

# Row 8 (Jaguares de Cordoba vs Pereira)
$ws.Range("G8").Value = 2.75
$ws.Range("I8").Value = 2.75
$ws.Range("L8").Value = 3.6
$ws.Range("W8").Value = 6.5
$ws.Range("X8").Value = 12
$ws.Range("Y8").Value = 12
$ws.Range("Z8").Value = 29
$ws.Range("AW8").Value = 4.5
$ws.Range("AX8").Value = 17
$ws.Range("AZ8").Value = 51
